$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'59.502.04"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.03%  '
$ws.Range('D3').Value = "'2.602.56"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.12%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = "'554.72"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.66%  '
$ws.Range('D6').Value = "'140.49"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.35%  '
$ws.Range('D7').Value = "'0.999"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = "'0.596"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.37%  '
$ws.Range('D9').Value = "'2.601.33"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.76%  '
$ws.Range('D10').Value = "'6.70"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.14%  '
$ws.Range('E11').Value = '  +1.08%  '
$ws.Range('D12').Value = "'0.161"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +6.78%  '
$ws.Range('D13').Value = "'0.359"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +5.21%  '
$ws.Range('D14').Value = "'3.042.00"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.58%  '
$ws.Range('D15').Value = "'59.501.08"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.89%  '
$ws.Range('D16').Value = "'23.06"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +5.06%  '
$ws.Range('E17').Value = '  +0.11%  '
$ws.Range('D18').Value = "'2.592.01"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.67%  '
$ws.Range('E19').Value = '  +0.87%  '
$ws.Range('D20').Value = "'340.53"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.75%  '
$ws.Range('D21').Value = "'10.45"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.18%  '
$ws.Range('D22').Value = "'6.59"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +6.95%  '
$ws.Range('D23').Value = "'0.998"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.21%  '
$ws.Range('D24').Value = "'0.491"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +10.51%  '
$ws.Range('D25').Value = "'62.19"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.27%  '
$ws.Range('D26').Value = "'0.997"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.12%  '
$ws.Range('D27').Value = "'0.159"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.47%  '
$ws.Range('D28').Value = "'7.47"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.47%  '
$ws.Range('D29').Value = "'0.0₃0772"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.83%  '
$ws.Range('E30').Value = '  +0.00%  '
$ws.Range('D31').Value = "'1.69"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.95%  '
$ws.Range('D32').Value = "'6.12"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.68%  '
$ws.Range('D33').Value = "'158.43"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.73%  '
$ws.Range('D34').Value = "'19.30"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.09%  '
$ws.Range('D35').Value = "'4.11"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +3.18%  '
$ws.Range('E36').Value = '  +3.98%  '
$ws.Range('D37').Value = "'1.16"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.20%  '
$ws.Range('D38').Value = "'37.54"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.99%  '
$ws.Range('D39').Value = "'1.49"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.97%  '
$ws.Range('D40').Value = "'0.840"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -4.11%  '
$ws.Range('D41').Value = "'3.67"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.47%  '
$ws.Range('D42').Value = "'291.35"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.64%  '
$ws.Range('D43').Value = "'136.31"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +9.33%  '
$ws.Range('D44').Value = "'0.999"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.07%  '
$ws.Range('D45').Value = "'0.0975"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.75%  '
$ws.Range('D46').Value = "'0.599"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.61%  '
$ws.Range('E47').Value = '  +0.33%  '
$ws.Range('D48').Value = "'0.0535"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.52%  '
$ws.Range('D49').Value = "'0.0236"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.33%  '

# Row 50/51 swap (Maker <-> InjectiveProtocol) with updated data
$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D50').Value = "'18.73"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.41%  '
$ws.Range('B51').Value = 'Maker'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D51').Value = "'1.964.28"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.99%  '
